# corregido bug de tipos de fixture
# - Bug #9 ("Tipo de Fixture") status moves from PENDIENTE to CORREGIDO.
# - The whole data table gets vertically centered (Align Middle).
# - The view scrolls down a bit and the active cell moves to F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark bug #9 ("Tipo de Fixture") as fixed.
$ws.Range("F11").Value = "CORREGIDO"

# Vertically center the whole table (header + data rows) - mirrors clicking
# the "Align Middle" button after selecting A1:F16.
$ws.Range("A1:F16").VerticalAlignment = -4108   # xlCenter

# Scroll the window down and move the selection, matching the author's
# on-screen state when the file was saved.
$ws.Range("F12").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
